$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (Changed) date column C for rows 2-7:
# 45184 (2023-09-15) -> 45185 (2023-09-16)
foreach ($row in 2..7) {
    $ws.Cells.Item($row, 3).Value = 45185
}
